$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Apply the "Alchemist" row styling (style used by rows 27-50, fill theme 7)
#        to rows 51-55 (previously blank "Mage" placeholder rows with a different fill),
#        columns A:C. Column D already shares the same style in both sections.
$ws.Range("A50:C50").Copy() | Out-Null
$ws.Range("A51:C51").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A50:C50").Copy() | Out-Null
$ws.Range("A52:C55").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# PasteSpecial(xlPasteFormats) only copies formatting, not values, so column A
# (Class = "Alchemist") still needs to be set explicitly for the rows that used
# to be blank "Mage" placeholders.
$ws.Cells.Item(51, 1).Value = "Alchemist"
$ws.Cells.Item(52, 1).Value = "Alchemist"
$ws.Cells.Item(53, 1).Value = "Alchemist"
$ws.Cells.Item(54, 1).Value = "Alchemist"
$ws.Cells.Item(55, 1).Value = "Alchemist"

# --- 2. Shift the existing "Healing_Potion / Poison_Potion / Mana_Potion / Antidote_Potion"
#        rows (currently B47:D50) down into B52:D55, preserving their Class (column A
#        already set to "Alchemist" by the paste above) and Robux Cost (column D = 20).
#        Use Value2 (not Value) so the scalar text/number is copied rather than a
#        COM variant wrapper.
$ws.Cells.Item(52, 2).Value = $ws.Cells.Item(47, 2).Value2
$ws.Cells.Item(52, 3).Value = $ws.Cells.Item(47, 3).Value2
$ws.Cells.Item(52, 4).Value = 20

$ws.Cells.Item(53, 2).Value = $ws.Cells.Item(48, 2).Value2
$ws.Cells.Item(53, 3).Value = $ws.Cells.Item(48, 3).Value2
$ws.Cells.Item(53, 4).Value = 20

$ws.Cells.Item(54, 2).Value = $ws.Cells.Item(49, 2).Value2
$ws.Cells.Item(54, 3).Value = $ws.Cells.Item(49, 3).Value2
$ws.Cells.Item(54, 4).Value = 20

$ws.Cells.Item(55, 2).Value = $ws.Cells.Item(50, 2).Value2
$ws.Cells.Item(55, 3).Value = $ws.Cells.Item(50, 3).Value2
$ws.Cells.Item(55, 4).Value = 20

# --- 3. Overwrite rows 47-51 with the five new "Potion_Throw" sound entries.
$ws.Cells.Item(47, 2).Value = "Potion_Throw_01"
$ws.Cells.Item(47, 3).Value = 0.41

$ws.Cells.Item(48, 2).Value = "Potion_Throw_02"
$ws.Cells.Item(48, 3).Value = 0.41

$ws.Cells.Item(49, 2).Value = "Potion_Throw_03"
$ws.Cells.Item(49, 3).Value = 0.34

$ws.Cells.Item(50, 2).Value = "Potion_Throw_04"
$ws.Cells.Item(50, 3).Value = 0.33

$ws.Cells.Item(51, 2).Value = "Potion_Throw_05"
$ws.Cells.Item(51, 3).Value = 0.43
$ws.Cells.Item(51, 4).Value = 20

# --- 4. Update sheet view selection to match the saved session state
#        (mirrors the <selection activeCell="D52" sqref="D52:D55"/> in the diff).
$ws.Range("D52:D55").Select() | Out-Null
